$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 6 new GPS-tracking rows (rows 496-501) for the J+3 training session
# dated 2025-09-16, mirroring the existing table's row layout (copy styles
# from a representative existing "J+3" row, then set the actual cell values).
# ---------------------------------------------------------------------------

# Use row 462 as the style template: it already carries the date format on
# column B (style index 1) and the centered "J+3" style on column D (style
# index 3), with plain/default formatting on every other column - exactly
# what rows 496-501 need.
$ws.Range("A462:V462").Copy() | Out-Null
$ws.Range("A496:V501").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 496
$ws.Range("A496").Value = "Entrainement"
$ws.Range("B496").Value = 45916
$ws.Range("C496").Value = "Global"
$ws.Range("D496").Value = "J+3"
$ws.Range("E496").Value = "Mattheo Haon"
$ws.Range("F496").Value = "right back"
$ws.Range("G496").Value = "01:48:28"
$ws.Range("H496").Value = 8.61
$ws.Range("I496").Value = 0.3
$ws.Range("J496").Value = 8.3000000000000007
$ws.Range("K496").Value = 0.24
$ws.Range("L496").Value = 0.06
$ws.Range("M496").Value = 0.01
$ws.Range("N496").Value = 0
$ws.Range("O496").Value = 1
$ws.Range("P496").Value = 4.7
$ws.Range("Q496").Value = 27.44
$ws.Range("R496").Value = 4.76
$ws.Range("S496").Value = 41
$ws.Range("T496").Value = 9
$ws.Range("U496").Value = 19
$ws.Range("V496").Value = 5
# Row 497
$ws.Range("A497").Value = "Entrainement"
$ws.Range("B497").Value = 45916
$ws.Range("C497").Value = "Global"
$ws.Range("D497").Value = "J+3"
$ws.Range("E497").Value = "Ilan Ihaddadene"
$ws.Range("F497").Value = "center midfield"
$ws.Range("G497").Value = "01:48:12"
$ws.Range("H497").Value = 9.11
$ws.Range("I497").Value = 0.22
$ws.Range("J497").Value = 8.89
$ws.Range("K497").Value = 0.21
$ws.Range("L497").Value = 0.01
$ws.Range("M497").Value = 0
$ws.Range("N497").Value = 0
$ws.Range("O497").Value = 0
$ws.Range("P497").Value = 4.97
$ws.Range("Q497").Value = 21.99
$ws.Range("R497").Value = 4.68
$ws.Range("S497").Value = 30
$ws.Range("T497").Value = 6
$ws.Range("U497").Value = 28
$ws.Range("V497").Value = 2
# Row 498
$ws.Range("A498").Value = "Entrainement"
$ws.Range("B498").Value = 45916
$ws.Range("C498").Value = "Global"
$ws.Range("D498").Value = "J+3"
$ws.Range("E498").Value = "Omar Benyounes"
$ws.Range("F498").Value = "center midfield"
$ws.Range("G498").Value = "01:48:03"
$ws.Range("H498").Value = 8.86
$ws.Range("I498").Value = 0.39
$ws.Range("J498").Value = 8.4600000000000009
$ws.Range("K498").Value = 0.36
$ws.Range("L498").Value = 0.04
$ws.Range("M498").Value = 0
$ws.Range("N498").Value = 0
$ws.Range("O498").Value = 0
$ws.Range("P498").Value = 4.83
$ws.Range("Q498").Value = 22.04
$ws.Range("R498").Value = 4.55
$ws.Range("S498").Value = 38
$ws.Range("T498").Value = 3
$ws.Range("U498").Value = 42
$ws.Range("V498").Value = 5
# Row 499
$ws.Range("A499").Value = "Entrainement"
$ws.Range("B499").Value = 45916
$ws.Range("C499").Value = "Global"
$ws.Range("D499").Value = "J+3"
$ws.Range("E499").Value = "Malik Boussaid"
$ws.Range("F499").Value = "right back"
$ws.Range("G499").Value = "01:47:09"
$ws.Range("H499").Value = 8.69
$ws.Range("I499").Value = 0.23
$ws.Range("J499").Value = 8.4600000000000009
$ws.Range("K499").Value = 0.23
$ws.Range("L499").Value = 0
$ws.Range("M499").Value = 0
$ws.Range("N499").Value = 0
$ws.Range("O499").Value = 0
$ws.Range("P499").Value = 4.75
$ws.Range("Q499").Value = 21.59
$ws.Range("R499").Value = 4.62
$ws.Range("S499").Value = 46
$ws.Range("T499").Value = 2
$ws.Range("U499").Value = 34
$ws.Range("V499").Value = 8
# Row 500
$ws.Range("A500").Value = "Entrainement"
$ws.Range("B500").Value = 45916
$ws.Range("C500").Value = "Global"
$ws.Range("D500").Value = "J+3"
$ws.Range("E500").Value = "Emmanuel Valey"
$ws.Range("F500").Value = "left forward"
$ws.Range("G500").Value = "01:41:21"
$ws.Range("H500").Value = 8.92
$ws.Range("I500").Value = 0.22
$ws.Range("J500").Value = 8.68
$ws.Range("K500").Value = 0.23
$ws.Range("L500").Value = 0
$ws.Range("M500").Value = 0
$ws.Range("N500").Value = 0
$ws.Range("O500").Value = 0
$ws.Range("P500").Value = 4.9400000000000004
$ws.Range("Q500").Value = 20.41
$ws.Range("R500").Value = 4.93
$ws.Range("S500").Value = 85
$ws.Range("T500").Value = 12
$ws.Range("U500").Value = 78
$ws.Range("V500").Value = 18
# Row 501
$ws.Range("A501").Value = "Entrainement"
$ws.Range("B501").Value = 45916
$ws.Range("C501").Value = "Global"
$ws.Range("D501").Value = "J+3"
$ws.Range("E501").Value = "Karahali Souaré"
$ws.Range("F501").Value = "right forward"
$ws.Range("G501").Value = "01:46:27"
$ws.Range("H501").Value = 8.2100000000000009
$ws.Range("I501").Value = 0.28000000000000003
$ws.Range("J501").Value = 7.92
$ws.Range("K501").Value = 0.25
$ws.Range("L501").Value = 0.04
$ws.Range("M501").Value = 0
$ws.Range("N501").Value = 0
$ws.Range("O501").Value = 0
$ws.Range("P501").Value = 4.18
$ws.Range("Q501").Value = 23.51
$ws.Range("R501").Value = 5.79
$ws.Range("S501").Value = 57
$ws.Range("T501").Value = 11
$ws.Range("U501").Value = 45
$ws.Range("V501").Value = 21


# Row 499's "Poste" cell (F499) carries the small distinguishing font style
# (style index 6) seen elsewhere in the sheet for this same value/label
# combination - copy that single cell's format across.
$ws.Range("F468").Copy() | Out-Null
$ws.Range("F499").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the selection/active cell left by the editor after entering the
# new rows.
$ws.Range("D504").Select() | Out-Null
